# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3205
$ws1.Range("F5").Value = 6876
$ws1.Range("F6").Value = 2043
$ws1.Range("F7").Value = 24
$ws1.Range("F12").Value = 22
$ws1.Range("F14").Value = 187

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3205
$ws4.Range("F6").Value = 6876
$ws4.Range("F7").Value = 2043
$ws4.Range("F8").Value = 24
$ws4.Range("F13").Value = 22
$ws4.Range("F15").Value = 187
